$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = 24
$ws.Range("A4").Value = 17
$ws.Range("A5").Value = 12
$ws.Range("E5").Value = 4.2
$ws.Range("F5").Value = 55
$ws.Range("A7").Value = 49
$ws.Range("A10").Value = 37
$ws.Range("A11").Value = 45
$ws.Range("E11").Value = 4.1
$ws.Range("F11").Value = 219
$ws.Range("A12").Value = 53
$ws.Range("A14").Value = 50
$ws.Range("A15").Value = 46
$ws.Range("A16").Value = 26
$ws.Range("A18").Value = 32
$ws.Range("A19").Value = 28
$ws.Range("A21").Value = 29
$ws.Range("A22").Value = 52
$ws.Range("F22").Value = 78
$ws.Range("A24").Value = 27
$ws.Range("A27").Value = 36
$ws.Range("A28").Value = 39
$ws.Range("C28").Value = 'Sport Clips Haircuts of Castleton Crossing'
$ws.Range("E28").Value = 4.4
$ws.Range("A29").Value = 33
$ws.Range("C29").Value = 'Sport Clips Haircuts of Geist Oaklandon'
$ws.Range("E29").Value = 4.6
$ws.Range("F29").Value = 127
$ws.Range("A30").Value = 58
$ws.Range("C30").Value = 'Sport Clips Haircuts of German Church Shops'
$ws.Range("E30").Value = 3.9
$ws.Range("F30").Value = 129
$ws.Range("A31").Value = 51
$ws.Range("C31").Value = 'Sport Clips Haircuts of Greenwood'
$ws.Range("E31").Value = 4.3
$ws.Range("F31").Value = 136
$ws.Range("A32").Value = 48
$ws.Range("C32").Value = 'Sport Clips Haircuts of Greenwood Plaza'
$ws.Range("F32").Value = 105
$ws.Range("A33").Value = 54
$ws.Range("C33").Value = 'Sport Clips Haircuts of Greenwood South'
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 136
$ws.Range("A34").Value = 42
$ws.Range("C34").Value = 'Sport Clips Haircuts of Greenwood Springs'
$ws.Range("E34").Value = 4.3
$ws.Range("F34").Value = 143
$ws.Range("A35").Value = 44
$ws.Range("C35").Value = 'Sport Clips Haircuts of North Keystone'
$ws.Range("E35").Value = 3.9
$ws.Range("F35").Value = 100
$ws.Range("A36").Value = 55
$ws.Range("C36").Value = 'Sport Clips Haircuts of Northfield Commons'
$ws.Range("E36").Value = 4.5
$ws.Range("F36").Value = 148
$ws.Range("A37").Value = 41
$ws.Range("C37").Value = 'Sport Clips Haircuts of Plainfield'
$ws.Range("E37").Value = 4.4
$ws.Range("F37").Value = 167
$ws.Range("A38").Value = 47
$ws.Range("C38").Value = 'Sport Clips Haircuts of Shadeland Place'
$ws.Range("E38").Value = 4.3
$ws.Range("F38").Value = 91
$ws.Range("A39").Value = 43
$ws.Range("C39").Value = 'Sport Clips Haircuts of South Bluff Crossing'
$ws.Range("F39").Value = 111
$ws.Range("A40").Value = 34
$ws.Range("C40").Value = 'Sport Graphics Inc'
$ws.Range("F40").Value = 10
$ws.Range("A41").Value = 38
$ws.Range("C41").Value = 'Sport Vision'
$ws.Range("E41").Value = 4.9
$ws.Range("F41").Value = 7
$ws.Range("A42").Value = 56
$ws.Range("C42").Value = 'Sport''n Image'
$ws.Range("E42").Value = 5
$ws.Range("F42").Value = 3
$ws.Range("A43").Value = 57
$ws.Range("C43").Value = 'Sport.ly'
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("A44").Value = 20
$ws.Range("C44").Value = 'Sports Corporation Inc'
$ws.Range("A45").Value = 10
$ws.Range("C45").Value = 'Sports Plus'
$ws.Range("E45").Value = 4.9
$ws.Range("F45").Value = 16
$ws.Range("A46").Value = 0
$ws.Range("C46").Value = 'Sports Spot'
$ws.Range("E46").Value = 3.3
$ws.Range("F46").Value = 13
$ws.Range("A47").Value = 35
$ws.Range("C47").Value = 'Sports Travel & Tickets'
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("A48").Value = 18
$ws.Range("A49").Value = 19
$ws.Range("A51").Value = 9
$ws.Range("A52").Value = 40
$ws.Range("A54").Value = 11
